# Update imputed values in result_data_KNN.xlsx (Update Name of Algo)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -13.376
$ws.Range("B9").Value = 5.241000000000001
$ws.Range("C9").Value = -10.623
$ws.Range("C11").Value = -12.818
$ws.Range("B13").Value = 6.005999999999999
$ws.Range("B16").Value = 5.442
$ws.Range("C16").Value = -13.157
$ws.Range("B18").Value = 5.302000000000001
$ws.Range("B20").Value = 6.861999999999999
$ws.Range("C23").Value = -13.346
$ws.Range("C24").Value = -12.743
$ws.Range("B26").Value = 5.474
$ws.Range("C26").Value = -12.522
$ws.Range("B27").Value = 6.196
$ws.Range("B29").Value = 5.347
$ws.Range("C34").Value = -12.285
$ws.Range("B35").Value = 7.645
$ws.Range("C35").Value = -12.373
$ws.Range("B36").Value = 8.044
$ws.Range("C44").Value = -12.869
$ws.Range("B45").Value = 5.755
$ws.Range("C48").Value = -11.531
$ws.Range("C49").Value = -13.403
$ws.Range("C52").Value = -12.09
$ws.Range("B55").Value = 4.839
$ws.Range("B57").Value = 5.730000000000001
$ws.Range("C66").Value = -11.496
$ws.Range("C67").Value = -10.941
$ws.Range("B69").Value = 5.362
$ws.Range("C73").Value = -12.091
$ws.Range("B76").Value = 6.11
$ws.Range("B78").Value = 8.439
$ws.Range("C78").Value = -11.552
$ws.Range("C80").Value = -11.992
$ws.Range("B82").Value = 5.486
$ws.Range("B83").Value = 6.027
$ws.Range("C91").Value = -13.485
$ws.Range("B93").Value = 5.357000000000001
$ws.Range("B97").Value = 5.781
$ws.Range("C97").Value = -10.736
$ws.Range("C99").Value = -11.833
$ws.Range("C104").Value = -13.226
